$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row-level updates derived from the diff (coin reorder/price refresh + Hora flag flip)
$rowUpdates = @(
    @{ Row = 2; D = '252.13'; G = '1' },
    @{ Row = 3; D = '22.11'; G = '1' },
    @{ Row = 4; D = '5.575'; G = '1' },
    @{ Row = 5; D = '0.05684'; G = '1' },
    @{ Row = 6; D = '6.482'; G = '1' },
    @{ Row = 7; D = '0.8055'; G = '1' },
    @{ Row = 8; D = '1.055'; G = '1' },
    @{ Row = 9; B = 'WazirX'; C = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D = '0.1438'; E = '8WazirXWRX'; G = '1' },
    @{ Row = 10; B = 'MandalaExchangeToken'; C = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D = '0.07293'; E = '9MandalaExchangeTokenMDX'; G = '1' },
    @{ Row = 11; B = 'LiechtensteinCryptoassetsExchange'; C = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D = '0.03169'; E = '10LiechtensteinCryptoassetsExchangeLCX'; G = '1' },
    @{ Row = 12; B = 'BitrueCoin'; C = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D = '0.02939'; E = '11BitrueCoinBTR'; G = '1' },
    @{ Row = 13; B = 'BitMartToken'; C = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D = '0.09265'; E = '12BitMartTokenBMX'; G = '1' },
    @{ Row = 14; B = 'BitForexToken'; C = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D = '0.001658'; E = '13BitForexTokenBF'; G = '1' },
    @{ Row = 15; B = 'MCDex'; C = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'; D = '3.209'; E = '14MCDexMCB'; G = '1' },
    @{ Row = 16; B = 'CoinExToken'; C = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'; D = '0.04785'; E = '15CoinExTokenCET'; G = '1' },
    @{ Row = 17; B = 'One'; C = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; D = '0.0005813'; E = '16OneONE'; G = '1' },
    @{ Row = 18; D = '0.006417'; G = '1' },
    @{ Row = 19; D = '0.005071'; G = '1' },
    @{ Row = 20; D = '0.001048'; G = '1' },
    @{ Row = 21; G = '1' },
    @{ Row = 22; D = '0.0003201'; G = '1' },
    @{ Row = 23; G = '1' },
    @{ Row = 24; D = '3.381'; G = '1' },
    @{ Row = 25; D = '2.112'; G = '1' },
    @{ Row = 26; D = '0.3269'; G = '1' },
    @{ Row = 27; D = '0.1277'; G = '1' },
    @{ Row = 28; G = '1' },
    @{ Row = 29; G = '1' },
    @{ Row = 30; G = '1' },
    @{ Row = 31; G = '1' },
    @{ Row = 32; G = '1' },
    @{ Row = 33; G = '1' },
    @{ Row = 34; G = '1' },
    @{ Row = 35; G = '1' },
    @{ Row = 36; G = '1' },
    @{ Row = 37; G = '1' },
    @{ Row = 38; G = '1' },
    @{ Row = 39; G = '1' },
    @{ Row = 40; D = '0.04158'; G = '1' },
    @{ Row = 41; B = 'BKEXToken'; C = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'; D = '0.1049'; E = '40BKEXTokenBKK'; G = '1' },
    @{ Row = 42; D = '0.002971'; G = '1' },
    @{ Row = 43; B = 'KickToken'; C = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; D = '0.006922'; E = '42KickTokenKICK'; G = '1' },
    @{ Row = 44; D = '0.009539'; G = '1' },
    @{ Row = 45; D = '0.00005645'; G = '1' },
    @{ Row = 46; G = '1' },
    @{ Row = 47; D = '0.7856'; G = '1' },
    @{ Row = 48; D = '0.01669'; G = '1' },
    @{ Row = 49; D = '0.00002101'; G = '1' },
    @{ Row = 50; G = '1' },
    @{ Row = 51; G = '1' }
)

foreach ($item in $rowUpdates) {
    $r = $item.Row
    if ($item.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $item.B }
    if ($item.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $item.C }
    if ($item.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = "'" + $item.D }
    if ($item.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $item.E }
    if ($item.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = "'" + $item.G }
}
